# Update the "K" column (column G) values for rows 2-33 on Sheet1.
# The save_data regeneration switched the stat source from "Strike#" to "K",
# so the per-game strikeout counts (s_vals) were recalculated and rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 4
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 2
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 3
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 3
    24 = 0
    25 = 0
    26 = 2
    27 = 0
    28 = 0
    29 = 0
    30 = 3
    31 = 0
    32 = 2
    33 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
